$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 3 (row 6), DF hours (column D): add hour -> 10 becomes 11
$ws.Range("D6").Value = 11
